# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Recomputed K values are written back into column G (header "K") for each
# outing row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    "G2"  = 1
    "G3"  = 1
    "G4"  = 1
    "G5"  = 1
    "G6"  = 3
    "G7"  = 3
    "G8"  = 3
    "G9"  = 2
    "G10" = 3
    "G11" = 1
    "G12" = 1
    "G13" = 1
    "G14" = 1
    "G15" = 1
    "G16" = 1
    "G17" = 2
    "G18" = 1
    "G19" = 2
    "G20" = 3
    "G21" = 1
    "G22" = 3
    "G23" = 3
    "G24" = 4
    "G25" = 1
    "G26" = 1
    "G27" = 2
    "G28" = 3
    "G29" = 2
    "G30" = 2
    "G31" = 3
    "G32" = 2
    "G33" = 2
    "G34" = 2
    "G35" = 3
    "G36" = 0
    "G37" = 1
    "G38" = 3
    "G39" = 2
    "G40" = 5
    "G41" = 2
    "G42" = 1
    "G43" = 2
    "G44" = 1
    "G45" = 3
    "G46" = 0
    "G47" = 2
    "G48" = 2
    "G49" = 3
    "G50" = 3
    "G51" = 2
    "G52" = 1
    "G53" = 2
    "G54" = 2
    "G55" = 3
    "G56" = 0
    "G57" = 2
    "G58" = 2
    "G59" = 2
    "G60" = 1
    "G61" = 2
    "G62" = 1
    "G63" = 3
    "G64" = 1
    "G65" = 2
    "G66" = 2
    "G67" = 2
}

foreach ($cellRef in $kValues.Keys) {
    $ws.Range($cellRef).Value = $kValues[$cellRef]
}

Write-Host "Updated $($kValues.Count) K values in column G"
